# Auto-generated: apply cached market-data value updates to Sephirot_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 8374.5
$ws.Range("I19").Value = 10499.333
$ws.Range("J19").Value = 2000
$ws.Range("K19").Value = 10499.333
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = -10324.333
$ws.Range("N19").Value = -2350
$ws.Range("H28").Value = 547.25
$ws.Range("I28").Value = 547.25
$ws.Range("K28").Value = 547.25
$ws.Range("M28").Value = -62.25
$ws.Range("H97").Value = 1245.5
$ws.Range("J97").Value = 1245.5
$ws.Range("L97").Value = 3736.5
$ws.Range("N97").Value = -4728.5
$ws.Range("H107").Value = 265.55554
$ws.Range("I107").Value = 266
$ws.Range("K107").Value = 266
$ws.Range("M107").Value = 1654
$ws.Range("H118").Value = 479.5
$ws.Range("J118").Value = 709
$ws.Range("L118").Value = 2127
$ws.Range("N118").Value = -5441
$ws.Range("H138").Value = 2357.8667
$ws.Range("J138").Value = 3834
$ws.Range("L138").Value = 11502
$ws.Range("N138").Value = -21782
$ws.Range("H141").Value = 1285.421
$ws.Range("I141").Value = 1285.421
$ws.Range("K141").Value = 3856.263
$ws.Range("M141").Value = 1323.737

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2080.75
$ws.Range("I32").Value = 1965.6757
$ws.Range("K32").Value = 1965.6757
$ws.Range("M32").Value = -1678.6757
$ws.Range("H61").Value = 2465.7334
$ws.Range("I61").Value = 1544.1818
$ws.Range("K61").Value = 1544.1818
$ws.Range("M61").Value = -1332.1818
$ws.Range("H74").Value = 1109.1111
$ws.Range("I74").Value = 1108.5
$ws.Range("K74").Value = 1108.5
$ws.Range("M74").Value = -234.5
$ws.Range("H77").Value = 1109.1111
$ws.Range("I77").Value = 1108.5
$ws.Range("K77").Value = 5542.5
$ws.Range("M77").Value = -1174.5
$ws.Range("H132").Value = 2641.923
$ws.Range("I132").Value = 1593.8889
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 4781.6667
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -2251.6667
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 2465.7334
$ws.Range("I136").Value = 1544.1818
$ws.Range("K136").Value = 4632.5454
$ws.Range("M136").Value = -2082.5454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3554.2856
$ws.Range("I86").Value = 2079.6667
$ws.Range("K86").Value = 2079.6667
$ws.Range("M86").Value = -956.6667000000002
$ws.Range("H89").Value = 3554.2856
$ws.Range("I89").Value = 2079.6667
$ws.Range("K89").Value = 10398.3335
$ws.Range("M89").Value = -4782.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6250
$ws.Range("I99").Value = 5500
$ws.Range("K99").Value = 5500
$ws.Range("M99").Value = -4002
$ws.Range("H126").Value = 6250
$ws.Range("I126").Value = 5500
$ws.Range("K126").Value = 16500
$ws.Range("M126").Value = -14030

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 58.75
$ws.Range("I12").Value = 47.75
$ws.Range("J12").Value = 69.75
$ws.Range("K12").Value = 143.25
$ws.Range("L12").Value = 209.25
$ws.Range("M12").Value = 29.75
$ws.Range("N12").Value = -555.25
$ws.Range("H64").Value = 2050
$ws.Range("I64").Value = 1062.5
$ws.Range("J64").Value = 6000
$ws.Range("K64").Value = 3187.5
$ws.Range("L64").Value = 18000
$ws.Range("M64").Value = -2917.5
$ws.Range("N64").Value = -18540
$ws.Range("H67").Value = 2050
$ws.Range("I67").Value = 1062.5
$ws.Range("J67").Value = 6000
$ws.Range("K67").Value = 3187.5
$ws.Range("L67").Value = 18000
$ws.Range("M67").Value = -2251.5
$ws.Range("N67").Value = -19872
$ws.Range("H75").Value = 403
$ws.Range("I75").Value = 299
$ws.Range("J75").Value = 465.4
$ws.Range("K75").Value = 897
$ws.Range("L75").Value = 1396.2
$ws.Range("M75").Value = 101
$ws.Range("N75").Value = -3392.2
$ws.Range("H78").Value = 403
$ws.Range("I78").Value = 299
$ws.Range("J78").Value = 465.4
$ws.Range("K78").Value = 2691
$ws.Range("L78").Value = 4188.599999999999
$ws.Range("M78").Value = 2301
$ws.Range("N78").Value = -14172.6
$ws.Range("H103").Value = 869.8333
$ws.Range("J103").Value = 953.8
$ws.Range("L103").Value = 2861.4
$ws.Range("N103").Value = -4619.4
$ws.Range("H121").Value = 708
$ws.Range("I121").Value = 509.33334
$ws.Range("J121").Value = 906.6667
$ws.Range("K121").Value = 1528.00002
$ws.Range("L121").Value = 2720.0001
$ws.Range("M121").Value = -218.0000199999999
$ws.Range("N121").Value = -5340.0001
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -900
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 5000
$ws.Range("J23").Value = 5000
$ws.Range("L23").Value = 5000
$ws.Range("N23").Value = -5446
$ws.Range("H57").Value = 22042.6
$ws.Range("I57").Value = 16697
$ws.Range("J57").Value = 30061
$ws.Range("K57").Value = 16697
$ws.Range("L57").Value = 30061
$ws.Range("M57").Value = -15877
$ws.Range("N57").Value = -31701
$ws.Range("H107").Value = 1883.7858
$ws.Range("I107").Value = 1340.5714
$ws.Range("J107").Value = 2427
$ws.Range("K107").Value = 1340.5714
$ws.Range("L107").Value = 2427
$ws.Range("M107").Value = 579.4286
$ws.Range("N107").Value = -6267
$ws.Range("H132").Value = 4449
$ws.Range("I132").Value = 4449
$ws.Range("K132").Value = 13347
$ws.Range("M132").Value = -10817

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 10005.333
$ws.Range("J20").Value = 10005.333
$ws.Range("L20").Value = 10005.333
$ws.Range("N20").Value = -10457.333
$ws.Range("H136").Value = 735.6
$ws.Range("I136").Value = 735.6
$ws.Range("K136").Value = 2206.8
$ws.Range("M136").Value = 343.1999999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 63376.332
$ws.Range("J75").Value = 75064.5
$ws.Range("L75").Value = 75064.5
$ws.Range("N75").Value = -76936.5
$ws.Range("H78").Value = 63376.332
$ws.Range("J78").Value = 75064.5
$ws.Range("L78").Value = 225193.5
$ws.Range("N78").Value = -234553.5
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("H136").Value = 1045.6666
$ws.Range("I136").Value = 1045.6666
$ws.Range("K136").Value = 3136.9998
$ws.Range("M136").Value = -586.9998000000001
$ws.Range("N126").ClearContents()
